# Updates the "LOS Galacticos" Yahoo Fantasy roster sheet:
#  - re-orders the existing 16 players
#  - drops Jose Alvarado from row 2 (he moves further down the list)
#  - adds a new player, Donte DiVincenzo, as a new 18th row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$players = @(
    @("Keyonte George","PG,SG","Utah Jazz"),
    @("Anthony Edwards","SG,SF","Minnesota Timberwolves"),
    @("Bradley Beal","PG,SG,SF","Phoenix Suns"),
    @("James Harden","PG,SG","LA Clippers"),
    @("Jayson Tatum","SF,PF","Boston Celtics"),
    @("Paul George","SG,SF,PF","Philadelphia 76ers"),
    @("RJ Barrett","SG,SF,PF","Toronto Raptors"),
    @("Draymond Green","PF,C","Golden State Warriors"),
    @("Zion Williamson","PF,C","New Orleans Pelicans"),
    @("Jaren Jackson Jr.","PF,C","Memphis Grizzlies"),
    @("Anfernee Simons","PG,SG","Portland Trail Blazers"),
    @("Ivica Zubac","C","LA Clippers"),
    @("Jose Alvarado","PG","New Orleans Pelicans"),
    @("Giannis Antetokounmpo","PF,C","Milwaukee Bucks"),
    @("Amen Thompson","SG,SF,PF","Houston Rockets"),
    @("Fred VanVleet","PG","Houston Rockets"),
    @("Donte DiVincenzo","PG,SG,SF","Minnesota Timberwolves")
)

$row = 2
foreach ($p in $players) {
    $ws.Cells.Item($row, 1).Value = $p[0]
    $ws.Cells.Item($row, 2).Value = $p[1]
    $ws.Cells.Item($row, 3).Value = $p[2]
    $row++
}
